# DBField added, app was tested, all works fine
#
# Adds a new worksheet "Лист4" (after "Лист3") containing a small
# "temperature / layer count" table with two helper-formula columns,
# and updates the selection/scroll state on "Лист1" left over from the
# editing session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new sheet at the end of the tab strip.
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "Лист4"

# ---------------------------------------------------------------------
# 2. Header row (order matters: it drives shared-string allocation
#    order -> Температура, Количество слойев, Сумма индексов, Средний индекс).
# ---------------------------------------------------------------------
$ws4.Range("A1").Value = "Температура"
$ws4.Range("B1").Value = "Количество слойев"
$ws4.Range("C1").Value = "Сумма индексов"
$ws4.Range("D1").Value = "Средний индекс"

# ---------------------------------------------------------------------
# 3. Data rows 2..42.
#    A: temperature, going 30 down to -10
#    B: layer count
#    C: sum of indices
#    D: =C/B               (average index)
#    E: running counter 1..41
#    F: =E/3
# ---------------------------------------------------------------------
$colA = @(30,29,28,27,26,25,24,23,22,21,20,19,18,17,16,15,14,13,12,11,10,9,8,7,6,5,4,3,2,1,0,-1,-2,-3,-4,-5,-6,-7,-8,-9,-10)
$colB = @(1,1,1,1,1,1,1,1,1,1,2,2,2,2,2,2,2,2,2,2,2,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3)
$colC = @(1,1,1,1,1,2,2,2,3,3,3,4,4,4,5,5,5,6,6,6,7,7,7,8,8,8,8,9,9,9,10,10,10,10,11,11,11,12,12,12,12)

for ($i = 0; $i -lt 41; $i++) {
    $r = $i + 2
    $ws4.Cells.Item($r, 1).Value = $colA[$i]
    $ws4.Cells.Item($r, 2).Value = $colB[$i]
    $ws4.Cells.Item($r, 3).Value = $colC[$i]
    $ws4.Cells.Item($r, 4).Formula = "=C" + $r + "/B" + $r
    $ws4.Cells.Item($r, 5).Value = $i + 1
    $ws4.Cells.Item($r, 6).Formula = "=E" + $r + "/3"
}

# ---------------------------------------------------------------------
# 4. Column widths roughly matching the authored layout.
# ---------------------------------------------------------------------
$ws4.Columns.Item(1).ColumnWidth = 13.140625
$ws4.Columns.Item(2).ColumnWidth = 18.7109375
$ws4.Columns.Item(3).ColumnWidth = 16.5703125
$ws4.Columns.Item(4).ColumnWidth = 16.5703125

# ---------------------------------------------------------------------
# 5. Restore the leftover selection/scroll state on "Лист1" (no longer
#    the active tab once "Лист4" exists).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Лист1")
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("G3").Select()

# ---------------------------------------------------------------------
# 6. Leave "Лист4" as the active sheet/tab with F2 selected.
# ---------------------------------------------------------------------
$ws4.Activate()
$ws4.Range("F2").Select()
